$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to be a bare A1:D20 grid of numbers (Valid/T/Z/p-value,
# one statistical test per row, no headers/labels). The edit moves that
# data one column to the right (B:E) and one row down (2:21), adds a
# header row (Valid/T/Z/p-value) and labels each row with the pair of
# metrics the test was run on in the new column A.

$labels = @(
    "CyclomaticComplexity(CC) & CyclomaticComplexity(CC)",
    "MaintainabilityIndex & MaintainabilityIndex",
    "NbUniqueOperands & NbUniqueOperands",
    "NbOperands & NbOperands",
    "NbOperands & EffortToImplement",
    "NbUniqueOperators & NbUniqueOperators",
    "NbUniqueOperators & EffortToImplement",
    "NbOperators & NbOperators",
    "NbOperators & EffortToImplement",
    "ProgramLength & ProgramLength",
    "ProgramLength & EffortToImplement",
    "VocabularySize & VocabularySize",
    "ProgramVolume & ProgramVolume",
    "DifficultyLevel & DifficultyLevel",
    "ProgramLevel & ProgramLevel",
    "EffortToImplement & NbOperands",
    "EffortToImplement & NbUniqueOperators",
    "EffortToImplement & NbOperators",
    "EffortToImplement & ProgramLength",
    "EffortToImplement & EffortToImplement"
)

# Snapshot the original A1:D20 data block first (Value2 returns plain
# numbers via COM; .Value misbehaves for this host's Range/Cells objects).
$orig = @()
for ($r = 1; $r -le 20; $r++) {
    $vals = @()
    for ($c = 1; $c -le 4; $c++) {
        $vals += $ws.Cells.Item($r, $c).Value2
    }
    $orig += , $vals
}

# Wipe the old block so no stray values are left behind once things move.
$ws.Range("A1:D20").ClearContents()

# New header row.
$ws.Cells.Item(1, 2).Value = "Valid"
$ws.Cells.Item(1, 3).Value = "T"
$ws.Cells.Item(1, 4).Value = "Z"
$ws.Cells.Item(1, 5).Value = "p-value"

# Re-seat each original row one row down, labelled in column A and with
# its four numbers shifted from A:D into B:E.
for ($i = 0; $i -lt $orig.Count; $i++) {
    $destRow = $i + 2
    $ws.Cells.Item($destRow, 1).Value = $labels[$i]
    $vals = $orig[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($destRow, $c + 2).Value = $vals[$c]
    }
}

# Column widths. The COM ColumnWidth setter snaps to Excel's internal
# pixel grid, so these inputs are chosen to land on the closest
# achievable stored width to the target (54.552101 / 6.596372 /
# 9.593605 / 9.593605 / 9.593605 characters).
$ws.Columns.Item(1).ColumnWidth = 53.65
$ws.Columns.Item(2).ColumnWidth = 5.83
$ws.Columns.Item(3).ColumnWidth = 8.83
$ws.Columns.Item(4).ColumnWidth = 8.83
$ws.Columns.Item(5).ColumnWidth = 8.83
